$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.428.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.444.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.76%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.526'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.109'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.82%  '
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("E11").Value = '  -2.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.346'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.37'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.891.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.97%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000172'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.419.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.443.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.70'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '320.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.78%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '639.03'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.17%  '
$ws.Range("E27").Value = '  -1.03%  '
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0941'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.49%  '
$ws.Range("E30").Value = '  -4.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.75'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.67%  '
$ws.Range("E32").Value = '  -3.57%  '
$ws.Range("E33").Value = '  -1.48%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("E35").Value = '  -4.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '151.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("E37").Value = '  -4.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.362'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.71%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₆0304'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '152.40'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.599'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0499'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0898'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.12%  '
